# Add 2022-Q3 data:
#  1. Insert a new "2022-Q3" worksheet, right before "2022-Q2", populated
#     with the quarter's fund-holding table (same layout as the other
#     quarterly sheets).
#  2. Update the "总计" (summary) sheet: insert a new first data row for
#     2022-Q3 (count=10, value=0.72), pushing the existing quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q3" sheet
# ---------------------------------------------------------------------
$target = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($target)
$newSheet.Name = "2022-Q3"

# Re-resolve the "2022-Q2" sheet by name now that the sheet collection has
# shifted (the handle obtained before the Add() call above tracks the
# newly-inserted sheet instead, not the original one).
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Reuse the existing "2022-Q2" sheet's layout/formatting (header row +
# styled column A) as a starting template, then overwrite every value.
$q2Sheet.Range("A1:H6").Copy($newSheet.Range("A1:H6"))
$q2Sheet.Range("A2:H2").Copy($newSheet.Range("A7:H7"))
$q2Sheet.Range("A2:H2").Copy($newSheet.Range("A8:H8"))
$q2Sheet.Range("A2:H2").Copy($newSheet.Range("A9:H9"))
$q2Sheet.Range("A2:H2").Copy($newSheet.Range("A10:H10"))
$q2Sheet.Range("A2:H2").Copy($newSheet.Range("A11:H11"))

# Columns B,D,E,F,G hold numeric-looking values that must stay text
# (fund codes with leading zeros, percentages kept as strings, etc.)
$newSheet.Range("B2:B11").NumberFormat = "@"
$newSheet.Range("D2:G11").NumberFormat = "@"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$rows = @(
    @("013623", "湘财周期轮动一年持有期混合",       "4.08", "93.99", "6.69", "0.2730", 5),
    @("007012", "湘财长顺混合A",                     "3.12", "94.24", "5.90", "0.1841", 8),
    @("008128", "湘财长源股票A",                     "1.32", "94.07", "6.84", "0.0903", 6),
    @("007013", "湘财长顺混合C",                     "1.29", "94.24", "5.90", "0.0761", 8),
    @("008129", "湘财长源股票C",                     "0.53", "94.07", "6.84", "0.0363", 6),
    @("010797", "长城优选回报六个月持有期混合A",     "2.83", "28.94", "1.03", "0.0291", 9),
    @("003456", "信澳新目标灵活配置混合",             "0.39", "94.17", "2.62", "0.0102", 1),
    @("014829", "诺德新能源汽车混合A",               "0.18", "81.13", "5.08", "0.0091", 7),
    @("014830", "诺德新能源汽车混合C",               "0.08", "81.13", "5.08", "0.0041", 7),
    @("010798", "长城优选回报六个月持有期混合C",     "0.30", "28.94", "1.03", "0.0031", 9)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push the existing data rows down by one, inheriting row 3's formatting
# for the freshly inserted row 2.
$summary.Rows.Item(2).Insert()
$summary.Range("A3:D3").Copy($summary.Range("A2:D2"))

# Column A is a 0-based rank matching each row's position, so every data
# row needs to be rewritten (not just the newly inserted one).
$summaryRows = @(
    @("2022-Q3", 10, 0.72),
    @("2022-Q2", 5, 0.24),
    @("2022-Q1", 7, 1.17),
    @("2021-Q4", 15, 4.49),
    @("2021-Q3", 11, 2.07),
    @("2021-Q2", 2, 1.1),
    @("2021-Q1", 1, 0.25),
    @("2020-Q4", 3, 0.11)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
}
